# Global_M2/Datasums/Long27_DataComp.xlsx update:
# Refreshed data-history columns (M2_Len/FX_Len counts and the
# M2/FX first/last date stamps) on Sheet1 for the rows whose source
# series picked up a new month of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - China / CNY
$ws.Range("G2").Value = 30074
$ws.Range("H2").Value = 45261

# Row 3 - United States / USD
$ws.Range("E3").Value = 30011
$ws.Range("F3").Value = 45200

# Row 4 - Euro Area / EUR
$ws.Range("E4").Value = 30011
$ws.Range("F4").Value = 45200
$ws.Range("G4").Value = 30074
$ws.Range("H4").Value = 45261

# Row 5 - Japan / JPY
$ws.Range("G5").Value = 30074
$ws.Range("H5").Value = 45261

# Row 6 - United Kingdom / GBP
$ws.Range("C6").Value = 443
$ws.Range("F6").Value = 45200
$ws.Range("G6").Value = 30074
$ws.Range("H6").Value = 45261

# Row 7 - South Korea / KRW
$ws.Range("G7").Value = 30074
$ws.Range("H7").Value = 45261

# Row 8 - Hong Kong / HKD
$ws.Range("D8").Value = 410
$ws.Range("E8").Value = 30011
$ws.Range("F8").Value = 45200
$ws.Range("H8").Value = 45261

# Row 9 - Australia / AUD
$ws.Range("E9").Value = 30011
$ws.Range("F9").Value = 45200
$ws.Range("G9").Value = 30074
$ws.Range("H9").Value = 45261

# Row 10 - Taiwan / TWD
$ws.Range("D10").Value = 483
$ws.Range("E10").Value = 30011
$ws.Range("F10").Value = 45200
$ws.Range("H10").Value = 45261

# Row 11 - Canada / CAD
$ws.Range("G11").Value = 30074
$ws.Range("H11").Value = 45261

# Row 12 - Russia / RUB
$ws.Range("C12").Value = 371
$ws.Range("D12").Value = 352
$ws.Range("F12").Value = 45200
$ws.Range("H12").Value = 45261

# Row 13 - Switzerland / CHF
$ws.Range("C13").Value = 467
$ws.Range("F13").Value = 45200
$ws.Range("G13").Value = 30074
$ws.Range("H13").Value = 45261

# Row 14 - Brazil / BRL
$ws.Range("D14").Value = 396
$ws.Range("H14").Value = 45261

# Row 15 - India / INR
$ws.Range("G15").Value = 30042
$ws.Range("H15").Value = 45261

# Row 16 - Mexico / MXN
$ws.Range("C16").Value = 455
$ws.Range("D16").Value = 410
$ws.Range("F16").Value = 45200
$ws.Range("H16").Value = 45261

# Row 17 - Saudi Arabia / SAR
$ws.Range("C17").Value = 370
$ws.Range("D17").Value = 394
$ws.Range("F17").Value = 45200
$ws.Range("H17").Value = 45261

# Row 18 - Singapore / SGD
$ws.Range("E18").Value = 30011
$ws.Range("F18").Value = 45200
$ws.Range("G18").Value = 30074
$ws.Range("H18").Value = 45261

# Row 19 - Indonesia / IDR
$ws.Range("D19").Value = 398
$ws.Range("E19").Value = 30011
$ws.Range("F19").Value = 45200
$ws.Range("H19").Value = 45261

# Row 20 - Malaysia / MYR
$ws.Range("C20").Value = 481
$ws.Range("F20").Value = 45200
$ws.Range("G20").Value = 30074
$ws.Range("H20").Value = 45261

# Row 21 - Norway / NOK
$ws.Range("E21").Value = 30011
$ws.Range("F21").Value = 45200
$ws.Range("G21").Value = 30074
$ws.Range("H21").Value = 45261

# Row 22 - Philippines / PHP
$ws.Range("D22").Value = 380
$ws.Range("H22").Value = 45261

# Row 23 - New Zealand / NZD
$ws.Range("D23").Value = 217
$ws.Range("E23").Value = 30011
$ws.Range("F23").Value = 45200
$ws.Range("H23").Value = 45261

# Row 24 - Denmark / DKK
$ws.Range("C24").Value = 394
$ws.Range("F24").Value = 45200
$ws.Range("G24").Value = 30074
$ws.Range("H24").Value = 45261

# Row 25 - South Africa / ZAR
$ws.Range("E25").Value = 30011
$ws.Range("F25").Value = 45200
$ws.Range("G25").Value = 30074
$ws.Range("H25").Value = 45261

# Row 26 - Chile / CLP
$ws.Range("D26").Value = 398
$ws.Range("H26").Value = 45261

# Row 27 - Colombia / COP
$ws.Range("D27").Value = 410
$ws.Range("H27").Value = 45261

# Row 28 - Kuwait / KWD
$ws.Range("C28").Value = 359
$ws.Range("D28").Value = 367
$ws.Range("F28").Value = 45200
$ws.Range("H28").Value = 45261
